# Trade #25 closed at 2026-02-18 00:13:00 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.86
$summary.Range("B4").Value = 0.96
$summary.Range("B5").Value = 0.36
$summary.Range("B6").Value = 53
$summary.Range("B7").Value = 30
$summary.Range("B9").Value = 56.6

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.86
$status.Range("D6").Value = 24
$status.Range("E6").Value = 0.05
$status.Range("F6").Value = -0.14
$status.Range("G6").Value = 58.33

# ---------------------------------------------------------------
# All Trades sheet - close trade #55 (row 56) and append trade #84 (row 85)
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G56").Value = 0.75
$allTrades.Range("H56").Value = "CLOSED"
$allTrades.Range("I56").Value = 7.1429
$allTrades.Range("J56").Value = 0.05
$allTrades.Range("K56").Value = 99.86
$allTrades.Range("L56").Value = "early_exit"
$allTrades.Range("M56").Value = 0.15

$allTrades.Range("A85").Value = 84
$allTrades.Range("B85").Value = "'2026-02-18"
$allTrades.Range("C85").Value = "00:12:54"
$allTrades.Range("D85").Value = "momentum"
$allTrades.Range("E85").Value = "DOWN"
$allTrades.Range("F85").Value = 0.7
$allTrades.Range("G85").NumberFormat = "General"
$allTrades.Range("H85").Value = "OPEN"
$allTrades.Range("I85").Value = 0
$allTrades.Range("J85").Value = 0
$allTrades.Range("K85").Value = 100
$allTrades.Range("L85").NumberFormat = "General"
$allTrades.Range("M85").Value = 0
$allTrades.Range("N85").Value = 0
$allTrades.Range("O85").Value = 0
$allTrades.Range("P85").Value = 0.9
$allTrades.Range("Q85").Value = "Downward momentum: -1.980% over 10 samples"

# ---------------------------------------------------------------
# momentum sheet - append new open trade #84 (row 15)
# ---------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A15").Value = 84
$momentum.Range("B15").Value = "'2026-02-18"
$momentum.Range("C15").Value = "00:12:54"
$momentum.Range("D15").Value = "momentum"
$momentum.Range("E15").Value = "DOWN"
$momentum.Range("F15").Value = 0.7
$momentum.Range("G15").NumberFormat = "General"
$momentum.Range("H15").Value = "OPEN"
$momentum.Range("I15").Value = 0
$momentum.Range("J15").Value = 0
$momentum.Range("K15").Value = 100
$momentum.Range("L15").Value = 0
$momentum.Range("M15").Value = 0
$momentum.Range("N15").Value = 0.9
$momentum.Range("O15").Value = "Downward momentum: -1.980% over 10 samples"
$momentum.Range("P15").NumberFormat = "General"
$momentum.Range("Q15").Value = 0

# ---------------------------------------------------------------
# MarketMaking sheet - close trade #55 (row 27)
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G27").Value = 0.75
$marketMaking.Range("H27").Value = "CLOSED"
$marketMaking.Range("I27").Value = 7.1429
$marketMaking.Range("J27").Value = 0.05
$marketMaking.Range("K27").Value = 99.86
$marketMaking.Range("P27").Value = "early_exit"
$marketMaking.Range("Q27").Value = 0.15
